$d = $word.ActiveDocument

$r0 = $d.Content.Find.Execute("Small-space growing fails because of constraints that don’t negotiate.", $true, $true, $false, $false, $false, $true, 1, $false, "Small-space growing fails because the penalties are immediate and compounding.", 2)
Write-Output "Replace 0: $r0"

$r1 = $d.Content.Find.Execute("First: **it’s physics, not vibes.** Containers are small soil volumes exposed to heat and wind. A short watering failure during a hot spell can depress production for the rest of that fruiting cycle.", $true, $true, $false, $false, $false, $true, 1, $false, "First: **it’s physics, not vibes.** Containers are small soil volumes exposed to heat and wind. One missed watering during a hot spell can collapse flowering, crack fruit, and set you back weeks.", 2)
Write-Output "Replace 1: $r1"

$r2 = $d.Content.Find.Execute("Second: **most crops don’t matter for “meaningful.”** A meaningful fraction of your food comes from a narrow set of crops that reliably convert light and water into calories. Most small-space advice steers you toward low-calorie wins that look busy and leave your pantry unchanged.", $true, $true, $false, $false, $false, $true, 1, $false, "Second: **most crops don’t matter for “meaningful.”** You don’t feel “meaningful” in the pantry unless you grow crops that carry calories or meal volume. Greens are valuable, but they don’t solve the food problem people think they’re solving.", 2)
Write-Output "Replace 2: $r2"

$r3 = $d.Content.Find.Execute("Third: **logistics beat technique.** Soil volume, irrigation reliability, pest exclusion, and replacement planting determine yield. Beginners overspend on seeds and underspend on containers, mix quality, and watering reliability.", $true, $true, $false, $false, $false, $true, 1, $false, "Third: **logistics beat technique.** The season is won by water reliability, soil volume, pest exclusion, and replacement planting. Most beginners learn this backwards—after spending money on seeds and losing plants to preventable stress.", 2)
Write-Output "Replace 3: $r3"

$r4 = $d.Content.Find.Execute("Fourth: **microclimate decides.** Balconies and patios amplify heat, wind, and reflected light. The same container can produce dramatically different results a few feet apart.", $true, $true, $false, $false, $false, $true, 1, $false, "Fourth: **microclimate decides.** Balconies and patios amplify heat, wind, and reflected light. I’ve seen the same crop thrive on one side of a balcony and stall on the other because the “hot corner” cooked the root zone.", 2)
Write-Output "Replace 4: $r4"

$r5 = $d.Content.Find.Execute("**Light, regular feeding beats rescue feeding:** neglect + “hero fertilizer” is a reliable way to invite pests and disorder.", $true, $true, $false, $false, $false, $true, 1, $false, "**Light, regular feeding beats rescue feeding:** neglect + “hero fertilizer” reliably produces pests and disorder.", 2)
Write-Output "Replace 5: $r5"

$r6 = $d.Content.Find.Execute("What matters: if you don’t have ≥6 hours direct sun (or equivalent light), “meaningful” will come from greens/beans, not calorie crops.", $true, $true, $false, $false, $false, $true, 1, $false, "What matters: if you don’t have ≥6 hours direct sun (or equivalent light), “meaningful” comes from greens/beans, not calorie crops.", 2)
Write-Output "Replace 6: $r6"

$r7 = $d.Content.Find.Execute("What matters: if you routinely miss waterings in summer, scale the system to your life or automate; willpower doesn’t fix heat.", $true, $true, $false, $false, $false, $true, 1, $false, "What matters: if you miss waterings in summer, automate or reduce the system; missed water is a yield tax you pay all season.", 2)
Write-Output "Replace 7: $r7"

$r8 = $d.Content.Find.Execute("What matters: you want to learn *one system*, not collect experiences.", $true, $true, $false, $false, $false, $true, 1, $false, "What matters: you want one stable system you can run through heat, pests, and travel—not a dozen experiments that all need attention.", 2)
Write-Output "Replace 8: $r8"

$r9 = $d.Content.Find.Execute("What matters: if pests are common where you live, “reactive” becomes a seasonal tax; exclusion is how you keep momentum.", $true, $true, $false, $false, $false, $true, 1, $false, "What matters: in pest-heavy areas, “reactive” becomes a weekly drain; exclusion is how you keep output steady.", 2)
Write-Output "Replace 9: $r9"

$r10 = $d.Content.Find.Execute("What matters at this stage is measurement. If weight isn’t tracked weekly, the system drifts into vibes.", $true, $true, $false, $false, $false, $true, 1, $false, "What matters at this stage is measurement because your memory lies. People remember the best harvest week and forget the empty weeks.", 2)
Write-Output "Replace 10: $r10"

$r11 = $d.Content.Find.Execute("What matters at this stage is honesty: light, time, travel, and budget. A plan that conflicts with your calendar will fail quietly.", $true, $true, $false, $false, $false, $true, 1, $false, "What matters at this stage is admitting what will break first—usually travel, heat, or watering. A plan that needs perfect attendance collapses the first time life gets busy.", 2)
Write-Output "Replace 11: $r11"

$r12 = $d.Content.Find.Execute("What matters at this stage is removing predictable failure points (soil volume, mulch, water reliability). Seeds are easy; stability is the work.", $true, $true, $false, $false, $false, $true, 1, $false, "What matters at this stage is preventing the classic container failure: a stressed root zone from heat + drying cycles. If stability isn’t built in, you end up “gardening” by emergency.", 2)
Write-Output "Replace 12: $r12"

$r13 = $d.Content.Find.Execute("What matters at this stage is consistent output: staggered planting and space reserved for replacements.", $true, $true, $false, $false, $false, $true, 1, $false, "What matters at this stage is keeping something harvestable most weeks, not maximizing one big harvest. Continuity is what makes the system feel real.", 2)
Write-Output "Replace 13: $r13"

$r14 = $d.Content.Find.Execute("What matters at this stage is rhythm: water, scout, harvest, replant. Rescue cycles consume time and produce less.", $true, $true, $false, $false, $false, $true, 1, $false, "What matters at this stage is rhythm. Once you fall into rescue cycles, you spend more time and harvest less.", 2)
Write-Output "Replace 14: $r14"

$r15 = $d.Content.Find.Execute("What matters at this stage is detachment: underperformers lose space; reliable performers earn space.", $true, $true, $false, $false, $false, $true, 1, $false, "What matters at this stage is detachment. Keeping sentimental underperformers is how small spaces stay unproductive.", 2)
Write-Output "Replace 15: $r15"

$r16 = $d.Content.Find.Execute("Small-space food is operations: measurement, reliability, replacement, sequencing. Structure orients you; consistency delivers the harvest.", $true, $true, $false, $false, $false, $true, 1, $false, "Small-space food is operations: measurement, reliability, replacement, sequencing. Structure orients you; consistency produces the harvest.", 2)
Write-Output "Replace 16: $r16"

$results = @($r0, $r1, $r2, $r3, $r4, $r5, $r6, $r7, $r8, $r9, $r10, $r11, $r12, $r13, $r14, $r15, $r16)
$failures = ($results | Where-Object { $_ -ne $true }).Count
if ($failures -gt 0) {
    throw "edit.ps1: $failures of $($results.Count) Find/Replace operations failed"
}
Write-Output "All $($results.Count) replacements applied successfully."
